# Refresh market-price derived columns (H:N) across the per-job Leve sheets.
# Values/structure below mirror the upstream scheduled-runner data pull.
$wb = $excel.ActiveWorkbook

# --- ALC: refresh price/profit columns ---
$ws = $wb.Worksheets.Item("ALC")
$updates = [ordered]@{
    "H125" = 2587.7827
    "I125" = 949.1111
    "J125" = 8487
    "K125" = 8541.999899999999
    "L125" = 76383
    "M125" = -6081.999899999999
    "N125" = -81303
    "H134" = 65321.934
    "I134" = 20709
    "J134" = 68508.57
    "K134" = 20709
    "L134" = 68508.57
    "M134" = -15639
    "N134" = -78648.57
    "H138" = 4687.29
    "I138" = 3076.9375
    "J138" = 4994.024
    "K138" = 9230.8125
    "L138" = 14982.072
    "M138" = -4090.8125
    "N138" = -25262.072
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- ARM: refresh price/profit columns ---
$ws = $wb.Worksheets.Item("ARM")
$updates = [ordered]@{
    "H32" = 18379.488
    "I32" = 16668.963
    "J32" = 35912.375
    "K32" = 16668.963
    "L32" = 35912.375
    "M32" = -16381.963
    "N32" = -36486.375
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- CUL: refresh price/profit columns ---
$ws = $wb.Worksheets.Item("CUL")
$updates = [ordered]@{
    "H12" = 32
    "I12" = 30
    "J12" = 32.727272
    "K12" = 90
    "L12" = 98.181816
    "M12" = 83
    "N12" = -444.181816
    "H68" = 1389.8247
    "I68" = 796.5625
    "J68" = 1970.9796
    "K68" = 2389.6875
    "L68" = 5912.9388
    "M68" = -1578.6875
    "N68" = -7534.9388
    "H71" = 1389.8247
    "I71" = 796.5625
    "J71" = 1970.9796
    "K71" = 7169.0625
    "L71" = 17738.8164
    "M71" = -3113.0625
    "N71" = -25850.8164
    "H87" = 15616.667
    "I87" = 7182
    "J87" = 21641.428
    "K87" = 21546
    "L87" = 64924.284
    "M87" = -20298
    "N87" = -67420.284
    "H90" = 15616.667
    "I90" = 7182
    "J90" = 21641.428
    "K90" = 64638
    "L90" = 194772.852
    "M90" = -58398
    "N90" = -207252.852
    "H107" = 339424
    "I107" = 641.25
    "J107" = 709005.2
    "K107" = 1923.75
    "L107" = 2127015.6
    "M107" = -3.75
    "N107" = -2130855.6
    "H120" = 13623.923
    "I120" = 6055.5
    "J120" = 15000
    "K120" = 18166.5
    "L120" = 45000
    "M120" = -13328.5
    "N120" = -54676
    "H121" = 1572.4166
    "I121" = 376.66666
    "J121" = 1971
    "K121" = 1129.99998
    "L121" = 5913
    "M121" = 180.0000199999999
    "N121" = -8533
    "H122" = 333809.88
    "I122" = 421.3125
    "J122" = 714825.4
    "K122" = 3791.8125
    "L122" = 6433428.600000001
    "M122" = -1341.8125
    "N122" = -6438328.600000001
    "H123" = 4375
    "I123" = 3000
    "J123" = 4571.4287
    "K123" = 9000
    "L123" = 13714.2861
    "M123" = -6550
    "N123" = -18614.2861
    "H124" = 3921.4285
    "I124" = 1225
    "J124" = 5000
    "K124" = 3675
    "L124" = 15000
    "M124" = 1235
    "N124" = -24820
    "H125" = 2975
    "I125" = 2800
    "J125" = 3000
    "K125" = 8400
    "L125" = 9000
    "M125" = -3480
    "N125" = -18840
    "H126" = 4538.4614
    "I126" = 1966.6666
    "J126" = 5310
    "K126" = 5899.9998
    "L126" = 15930
    "M126" = -959.9997999999996
    "N126" = -25810
    "H127" = 1895.8462
    "I127" = 1000
    "J127" = 1970.5
    "K127" = 3000
    "L127" = 5911.5
    "M127" = 1960
    "N127" = -15831.5
    "H128" = 497900
    "I128" = 497900
    "J128" = 0
    "K128" = 1493700
    "L128" = 0
    "M128" = -1488720
    "H129" = 57148.668
    "I129" = 1141.2858
    "J129" = 92789.73
    "K129" = 3423.8574
    "L129" = 278369.19
    "M129" = 1576.1426
    "N129" = -288369.19
    "H130" = 2993.3333
    "I130" = 1192
    "J130" = 12000
    "K130" = 3576
    "L130" = 36000
    "M130" = 1444
    "N130" = -46040
    "H131" = 15184220
    "I131" = 62625384
    "J131" = 3047.2
    "K131" = 187876152
    "L131" = 9141.599999999999
    "M131" = -187871112
    "N131" = -19221.6
    "H132" = 900
    "I132" = 900
    "J132" = 0
    "K132" = 8100
    "L132" = 0
    "M132" = -5570
    "H133" = 1250
    "I133" = 1250
    "J133" = 0
    "K133" = 3750
    "L133" = 0
    "M133" = 1310
    "H134" = 6239.9287
    "I134" = 4979.857
    "J134" = 7500
    "K134" = 14939.571
    "L134" = 22500
    "M134" = -9869.571
    "N134" = -32640
    "H136" = 3998.9
    "I136" = 2522.25
    "J136" = 4983.3335
    "K136" = 7566.75
    "L136" = 14950.0005
    "M136" = -2466.75
    "N136" = -25150.0005
    "H137" = 36938.906
    "I137" = 3075
    "J137" = 48226.875
    "K137" = 9225
    "L137" = 144680.625
    "M137" = -4125
    "N137" = -154880.625
    "H138" = 2610
    "I138" = 1900
    "J138" = 4598
    "K138" = 5700
    "L138" = 13794
    "M138" = -560
    "N138" = -24074
    "H139" = 1224.36
    "I139" = 863.7619
    "J139" = 3117.5
    "K139" = 2591.2857
    "L139" = 9352.5
    "M139" = 2548.7143
    "N139" = -19632.5
    "H140" = 2947.8572
    "I140" = 2947.8572
    "J140" = 0
    "K140" = 8843.5716
    "L140" = 0
    "M140" = -3663.571599999999
    "H141" = 6124.4443
    "I141" = 1665
    "J141" = 9692
    "K141" = 4995
    "L141" = 29076
    "M141" = 185
    "N141" = -39436
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- GSM: refresh price/profit columns ---
$ws = $wb.Worksheets.Item("GSM")
$updates = [ordered]@{
    "H70" = 4326.654
    "I70" = 4086.0454
    "J70" = 5650
    "K70" = 4086.0454
    "L70" = 5650
    "M70" = -3816.0454
    "N70" = -6190
    "H73" = 4326.654
    "I73" = 4086.0454
    "J73" = 5650
    "K73" = 4086.0454
    "L73" = 5650
    "M73" = -3150.0454
    "N73" = -7522
    "H80" = 4346.154
    "I80" = 4416.6665
    "K80" = 4416.6665
    "M80" = -3418.6665
    "H83" = 4346.154
    "I83" = 4416.6665
    "K83" = 22083.3325
    "M83" = -17091.3325
}
foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}

# --- LTW: clear stale price/profit columns (leve no longer priced) ---
$ws = $wb.Worksheets.Item("LTW")
$clears = @(
    "H124", "I124", "J124", "K124", "L124", "N124", "H125", "I125",
    "J125", "K125", "L125", "N125", "H127", "I127", "J127", "K127",
    "L127", "N127", "H128", "I128", "J128", "K128", "L128", "N128",
    "H129", "I129", "J129", "K129", "L129", "N129", "H130", "I130",
    "J130", "K130", "L130", "H131", "I131", "J131", "K131", "L131",
    "H132", "I132", "J132", "K132", "L132", "M132", "N132", "H133",
    "I133", "J133", "K133", "L133", "N133", "H134", "I134", "J134",
    "K134", "L134", "N134", "H135", "I135", "J135", "K135", "L135",
    "M135", "N135", "H136", "I136", "J136", "K136", "L136", "N136",
    "H137", "I137", "J137", "K137", "L137", "N137", "H138", "I138",
    "J138", "K138", "L138", "N138", "H139", "I139", "J139", "K139",
    "L139", "N139", "H140", "I140", "J140", "K140", "L140", "N140",
    "H141", "I141", "J141", "K141", "L141", "N141"
)
foreach ($addr in $clears) {
    $ws.Range($addr).ClearContents()
}

Write-Output "Applied Leve price/profit refresh across ALC, ARM, CUL, GSM, LTW"
